$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.940.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +4.91%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.514.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.77%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'592.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.89%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'168.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +6.58%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.516.14"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.78%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.577"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.45%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.32%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +5.49%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.439"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +4.10%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.120.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.81%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.06%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +3.94%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +4.05%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'66.895.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +4.74%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.522.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.72%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.94%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'14.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +3.15%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'394.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.44%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'7.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.01%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'73.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.27%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +9.91%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.997"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.31%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +3.33%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'10.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +5.05%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +2.13%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.06%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'6.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +5.53%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +5.83%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.10%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'23.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.01%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'7.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +7.24%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.00%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +5.70%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'161.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.20%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.901"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +6.71%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +5.42%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0750"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.90%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'4.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +7.14%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +1.98%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'6.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.81%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.834.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.80%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'43.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.02%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'InjectiveProtocol"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'26.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.23%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'dogwifhat"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'2.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +6.40%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +3.57%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'352.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +5.89%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +4.60%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'33.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +11.28%  "
$ws.Range("E51").Style = "Normal"
